$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2020 column (N) mirroring the existing 2019 column (M) formatting.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2.1

$excel.CutCopyMode = 0

# Match the saved selection state (cell N9 selected on sheet 1).
$ws.Range("N9").Select() | Out-Null
